# Updated symbol list (price + 1h volume/change %) per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "270.14"
Set-TextValue "E2" "3.25%"
Set-TextValue "E3" "-1.77%"
Set-TextValue "D4" "4.699"
Set-TextValue "E4" "-0.16%"
Set-TextValue "D5" "0.06103"
Set-TextValue "E5" "-1.63%"
Set-TextValue "D6" "6.738"
Set-TextValue "E6" "0.32%"
Set-TextValue "D7" "0.8585"
Set-TextValue "E7" "1.02%"
Set-TextValue "D8" "0.8930"
Set-TextValue "E8" "-2.52%"
Set-TextValue "D9" "0.1418"
Set-TextValue "E9" "0.81%"
Set-TextValue "D10" "0.05034"
Set-TextValue "E10" "8.75%"
Set-TextValue "D11" "0.07097"
Set-TextValue "E11" "0.10%"
Set-TextValue "D12" "0.03169"
Set-TextValue "E12" "0.78%"
Set-TextValue "D13" "0.09031"
Set-TextValue "E13" "-0.24%"
Set-TextValue "D14" "0.001528"
Set-TextValue "E14" "-0.18%"
Set-TextValue "D15" "0.0006063"
Set-TextValue "E15" "-1.64%"
Set-TextValue "D16" "0.006089"
Set-TextValue "E16" "-0.67%"
Set-TextValue "E17" "-0.15%"
Set-TextValue "D18" "3.167"
Set-TextValue "E18" "-0.01%"
Set-TextValue "D19" "2.243"
Set-TextValue "E19" "2.94%"
Set-TextValue "E20" "-0.62%"
Set-TextValue "E21" "-0.78%"
Set-TextValue "D22" "3.842"
Set-TextValue "E22" "-5.83%"
Set-TextValue "D23" "0.04241"
Set-TextValue "E23" "-0.19%"
Set-TextValue "E24" "-2.02%"
Set-TextValue "D25" "0.004152"
Set-TextValue "E25" "9.20%"
Set-TextValue "D26" "0.0001200"
Set-TextValue "E26" "-0.01%"
Set-TextValue "D27" "0.0001680"
Set-TextValue "E27" "4.93%"
Set-TextValue "D40" "0.03955"
Set-TextValue "E40" "0.98%"
Set-TextValue "D41" "0.1117"
Set-TextValue "E41" "0.42%"
Set-TextValue "D42" "0.004186"
Set-TextValue "E42" "1.23%"
Set-TextValue "D43" "0.002010"
Set-TextValue "E43" "-7.96%"
Set-TextValue "D44" "0.01257"
Set-TextValue "E44" "-9.66%"
Set-TextValue "D45" "0.00005132"
Set-TextValue "E45" "-0.80%"
Set-TextValue "E46" "-0.10%"
Set-TextValue "D48" "0.2081"
Set-TextValue "E48" "24.89%"
Set-TextValue "D49" "0.00002099"
Set-TextValue "E49" "-0.10%"
Set-TextValue "D50" "0.0001999"
Set-TextValue "E50" "-0.10%"
